$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Grow the bottom "Top Campaigns" panel (shape id=36) ---
$panel = $s.Shapes.Item(9)
$panel.Height = 153.0708695417323

# --- 2. Add three new thin rounded-rectangle "row" bars under that panel ---
# Use the panel itself as a style/template source via Duplicate() so the
# new shapes inherit the same p:style / txBody structure, then reposition,
# resize and recolor them to match the target rows.

# id=3 -> "Rectangle: Rounded Corners 2"
$dup1 = $panel.Duplicate()
$row1 = $dup1.Item(1)
$row1.Name = "Rectangle: Rounded Corners 2"
$row1.Left = 90.23614173228347
$row1.Top = 420.0696062992126
$row1.Width = 779.5275590551181
$row1.Height = 21.165826771653542
$row1.Fill.ForeColor.RGB = 0xFBFBFB

# id=4 -> created then removed (burns an id/name so the next shape becomes 5/"...4")
$dupTmp = $panel.Duplicate()
$rowTmp = $dupTmp.Item(1)
$rowTmp.Delete()

# id=5 -> "Rectangle: Rounded Corners 4"
$dup2 = $panel.Duplicate()
$row2 = $dup2.Item(1)
$row2.Name = "Rectangle: Rounded Corners 4"
$row2.Left = 89.56787491574804
$row2.Top = 460.49772653543306
$row2.Width = 779.5275590551181
$row2.Height = 21.165826771653542
$row2.Fill.ForeColor.RGB = 0xFBFBFB

# id=6 -> "Rectangle: Rounded Corners 5"
$dup3 = $panel.Duplicate()
$row3 = $dup3.Item(1)
$row3.Name = "Rectangle: Rounded Corners 5"
$row3.Left = 89.56787491574804
$row3.Top = 500.92582707165354
$row3.Width = 779.5275590551181
$row3.Height = 21.259842919685042
$row3.Fill.ForeColor.RGB = 0xFBFBFB
